$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("L:L").Insert()

$ws.Range("K5").Value = "arg1"
$ws.Range("L5").Value = "arg2"

for ($row = 6; $row -le 20; $row++) {
  $ws.Range("K$row").Value = "1 or True"
  $ws.Range("L$row").Value = "0 or False"
}
Write-Output "done"
